$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $Addr, $Val)
    $r = $Sheet.Range($Addr)
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '25.640.55'
$ws.Range('E2').Value = '  -4.34%  '
$ws.Range('D3').Value = '1.811.85'
$ws.Range('E3').Value = '  -3.20%  '
Set-TextValue $ws 'D4' '1.002'
$ws.Range('E4').Value = '  +0.12%  '
Set-TextValue $ws 'D5' '278.13'
$ws.Range('E5').Value = '  -7.59%  '
Set-TextValue $ws 'D6' '1.003'
$ws.Range('E6').Value = '  +0.16%  '
Set-TextValue $ws 'D7' '0.5083'
$ws.Range('E7').Value = '  -4.96%  '
Set-TextValue $ws 'D8' '0.3520'
$ws.Range('E8').Value = '  -5.92%  '
Set-TextValue $ws 'D9' '44.34'
$ws.Range('E9').Value = '  -2.28%  '
Set-TextValue $ws 'D10' '0.06683'
$ws.Range('E10').Value = '  -7.03%  '
Set-TextValue $ws 'D11' '19.83'
$ws.Range('E11').Value = '  -8.21%  '
Set-TextValue $ws 'D12' '0.8191'
$ws.Range('E12').Value = '  -7.88%  '
Set-TextValue $ws 'D13' '0.07868'
$ws.Range('E13').Value = '  -3.74%  '
$ws.Range('D14').Value = '1.816.76'
$ws.Range('E14').Value = '  -2.82%  '
Set-TextValue $ws 'D15' '5.061'
$ws.Range('E15').Value = '  -4.77%  '
Set-TextValue $ws 'D16' '87.33'
$ws.Range('E16').Value = '  -6.05%  '
$ws.Range('E17').Value = '  -0.05%  '
Set-TextValue $ws 'D18' '14.05'
$ws.Range('E18').Value = '  -5.39%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws 'D19' '1.001'
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws 'D20' '0.000008005'
$ws.Range('E20').Value = '  -5.98%  '
$ws.Range('D21').Value = '25.742.40'
$ws.Range('E21').Value = '  -4.09%  '
Set-TextValue $ws 'D22' '4.737'
$ws.Range('E22').Value = '  -4.89%  '
Set-TextValue $ws 'D23' '9.982'
$ws.Range('E23').Value = '  -6.04%  '
Set-TextValue $ws 'D24' '6.109'
$ws.Range('E24').Value = '  -4.16%  '
Set-TextValue $ws 'D25' '2.242'
$ws.Range('E25').Value = '  -2.71%  '
Set-TextValue $ws 'D26' '142.46'
$ws.Range('E26').Value = '  -2.40%  '
Set-TextValue $ws 'D27' '1.663'
$ws.Range('E27').Value = '  -4.21%  '
Set-TextValue $ws 'D28' '17.10'
$ws.Range('E28').Value = '  -5.35%  '
Set-TextValue $ws 'D29' '108.92'
$ws.Range('E29').Value = '  -4.31%  '
Set-TextValue $ws 'D30' '4.304'
$ws.Range('E30').Value = '  -8.76%  '
Set-TextValue $ws 'D31' '4.219'
$ws.Range('E31').Value = '  -8.67%  '
Set-TextValue $ws 'D32' '0.08727'
$ws.Range('E32').Value = '  -4.70%  '
Set-TextValue $ws 'D33' '0.04849'
$ws.Range('E33').Value = '  -3.44%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws 'D34' '2.900'
$ws.Range('E34').Value = '  -1.32%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws 'D35' '0.7242'
$ws.Range('E35').Value = '  -10.06%  '
Set-TextValue $ws 'D36' '1.128'
$ws.Range('E36').Value = '  -4.17%  '
Set-TextValue $ws 'D37' '3.161'
$ws.Range('E37').Value = '  -1.00%  '
Set-TextValue $ws 'D38' '2.359'
$ws.Range('E38').Value = '  -12.09%  '
Set-TextValue $ws 'D39' '0.01847'
$ws.Range('E39').Value = '  -5.32%  '
Set-TextValue $ws 'D40' '0.5137'
$ws.Range('E40').Value = '  -15.86%  '
Set-TextValue $ws 'D41' '0.9664'
$ws.Range('E41').Value = '  -9.13%  '
Set-TextValue $ws 'D42' '114.81'
$ws.Range('E42').Value = '  +0.15%  '
Set-TextValue $ws 'D43' '6.234'
$ws.Range('E43').Value = '  -4.79%  '
Set-TextValue $ws 'D44' '7.981'
$ws.Range('E44').Value = '  -9.05%  '
Set-TextValue $ws 'D45' '1.000'
$ws.Range('E45').Value = '  -0.04%  '
Set-TextValue $ws 'D46' '0.4525'
$ws.Range('E46').Value = '  -13.32%  '
Set-TextValue $ws 'D47' '0.1365'
$ws.Range('E47').Value = '  -8.51%  '
$ws.Range('E48').Value = '  -3.13%  '
Set-TextValue $ws 'D49' '9.232'
$ws.Range('E49').Value = '  -7.31%  '
Set-TextValue $ws 'D50' '1.497'
$ws.Range('E50').Value = '  -9.42%  '
Set-TextValue $ws 'D51' '0.05815'
$ws.Range('E51').Value = '  -3.99%  '
